# Insert a new weekly price-report row for "Ramas de apio" just above the
# current row 6. This shifts the former rows 6-12 down to rows 7-13
# (preserving all of their data), and populates the new row 6 with the
# latest report: same as the former last entry (old row 12 / Americana (o))
# but for variety "Sin especificar" at a newer date.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push existing rows 6.. down by one, duplicating row 6's formatting
# (matches Excel's native Insert-with-shift-down behaviour).
$ws.Rows.Item(6).Insert()

# Fill in the newly inserted row 6 with the new observation.
$ws.Cells.Item(6, 1).Value = 1
$ws.Cells.Item(6, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(6, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(6, 4).Value = 44575
$ws.Cells.Item(6, 5).Value = 15
$ws.Cells.Item(6, 6).Value = 100112017
$ws.Cells.Item(6, 7).Value = "Ramas de apio"
$ws.Cells.Item(6, 8).Value = "Sin especificar"
$ws.Cells.Item(6, 9).Value = "Primera"
$ws.Cells.Item(6, 10).Value = 160
$ws.Cells.Item(6, 11).Value = 6500
$ws.Cells.Item(6, 12).Value = 7000
$ws.Cells.Item(6, 13).Value = 6750
$ws.Cells.Item(6, 14).Value = "`$/atado 7 kilos"
$ws.Cells.Item(6, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(6, 16).Value = 6750
$ws.Cells.Item(6, 17).Value = 1
$ws.Cells.Item(6, 18).Value = "Hortaliza"
